$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 318, shifting existing rows 318..418 down to 319..419
$ws.Rows.Item(318).EntireRow.Insert()

# Populate the newly inserted row 318 with the new record
$ws.Range("A318").Value = 5
$ws.Range("B318").Value = "Macroferia Regional de Talca"
$ws.Range("C318").Value = "Maule"
$ws.Range("D318").Value = 44524
$ws.Range("D318").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E318").Value = 7
$ws.Range("F318").Value = 100112004
$ws.Range("G318").Value = "Cebolla"
$ws.Range("H318").Value = "Sin especificar"
$ws.Range("I318").Value = "1a nueva(o)"
$ws.Range("J318").Value = 50000
$ws.Range("K318").Value = 1000
$ws.Range("L318").Value = 1000
$ws.Range("M318").Value = 1000
$ws.Range("N318").Value = "$/paquete 10 unidades (volumen en unidades)"
$ws.Range("O318").Value = "Región de O'Higgins"
$ws.Range("P318").Value = 100
$ws.Range("Q318").Value = 10
$ws.Range("R318").Value = "Hortaliza"
